$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment in I15
$ws.Range("I15").Value = "Extra 19.68 hours carry over"

# Row 36: Forrest Cordova -> Jamie, hours 40 -> 20
$ws.Range("A36").Value = "Jamie"
$ws.Range("G36").Value = 20

# Row 37: Jamie -> Ayooluwa
$ws.Range("A37").Value = "Ayooluwa"

# Row 38: Ayooluwa -> Akshat
$ws.Range("A38").Value = "Akshat"

# Row 39: Akshat -> DJ
$ws.Range("A39").Value = "DJ"

# Remove row 40 entirely (was DJ row, now merged/removed)
$ws.Rows.Item(40).Delete()
